$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new rows above the current row 2 (the KOSS trade), pushing the
# existing trades (KOSS..ZAPP) down to rows 12-18.
$ws.Rows("2:11").Insert()

# Column A (dates) should use the same date-formatted style as the rest
# of the table; copy that formatting down from row 12 (the original row 2).
$ws.Range("A12").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Fill in the 10 new trade rows (rows 2-11) ----

# Row 2: GDHG
$ws.Range("A2").Value = 45474
$ws.Range("B2").Value = "GDHG"
$ws.Range("C2").Value = 0.1895
$ws.Range("D2").Value = 0.17799999999999999
$ws.Range("E2").Value = 500

# Row 3: LPTV
$ws.Range("A3").Value = 45474
$ws.Range("B3").Value = "LPTV"
$ws.Range("C3").Value = 0.15060000000000001
$ws.Range("D3").Value = 0.14000000000000001
$ws.Range("E3").Value = 500

# Row 4: SNTI
$ws.Range("A4").Value = 45474
$ws.Range("B4").Value = "SNTI"
$ws.Range("C4").Value = 0.48509999999999998
$ws.Range("D4").Value = 0.52
$ws.Range("E4").Value = 500

# Row 5: SNTI
$ws.Range("A5").Value = 45474
$ws.Range("B5").Value = "SNTI"
$ws.Range("C5").Value = 0.47
$ws.Range("D5").Value = 0.52
$ws.Range("E5").Value = 400

# Row 6: LUCY
$ws.Range("A6").Value = 45474
$ws.Range("B6").Value = "LUCY"
$ws.Range("C6").Value = 0.51649999999999996
$ws.Range("D6").Value = 0.49209999999999998
$ws.Range("E6").Value = 500
$ws.Rows("6").RowHeight = 23

# Row 7: GOVX
$ws.Range("A7").Value = 45474
$ws.Range("B7").Value = "GOVX"
$ws.Range("C7").Value = 3.61
$ws.Range("D7").Value = 3.55
$ws.Range("E7").Value = 65
$ws.Rows("7").RowHeight = 23

# Row 8: TSLA
$ws.Range("A8").Value = 45475
$ws.Range("B8").Value = "TSLA"
$ws.Range("C8").Value = 210.82
$ws.Range("D8").Value = 226.46
$ws.Range("E8").Value = 2
$ws.Rows("8").RowHeight = 23

# Row 9: OPTT
$ws.Range("A9").Value = 45475
$ws.Range("B9").Value = "OPTT"
$ws.Range("C9").Value = 0.45500000000000002
$ws.Range("D9").Value = 0.54
$ws.Range("E9").Value = 100
$ws.Rows("9").RowHeight = 23

# Row 10: MBIO
$ws.Range("A10").Value = 45475
$ws.Range("B10").Value = "MBIO"
$ws.Range("C10").Value = 0.53500000000000003
$ws.Range("D10").Value = 0.56599999999999995
$ws.Range("E10").Value = 100
$ws.Rows("10").RowHeight = 23

# Row 11: ANVS
$ws.Range("A11").Value = 45475
$ws.Range("B11").Value = "ANVS"
$ws.Range("C11").Value = 6.47
$ws.Range("D11").Value = 10.71
$ws.Range("E11").Value = 1
$ws.Rows("11").RowHeight = 23

# The old first trade row (KOSS, now row 12) has its "Stock name" cell
# restyled to match the plain style used by the rest of the newly-entered
# rows (it used to share the date-number-format style with column A).
$ws.Range("B2").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Update the active selection to match where the user ended up ----
$ws.Range("G16").Select() | Out-Null
